# "changed around pflow and qflow. still not quite working, need to fix formula."
#
# - sheet "line_imp" (2nd sheet): a new row is inserted above the old row 4,
#   pushing the old rows 4/5/6 down to 5/6/7. The (new) row 3 has its B/C
#   formulas replaced with plain numeric values and its E formula simplified,
#   and the freshly-inserted row 4 is filled with a duplicate of that same
#   row-3 data. Row 5 (old row 4) gets its bus number corrected back to 14.
# - sheet "initial" (1st sheet): only the selected cell changes (no data
#   edits).
#
# Selections are applied last, sheet2 before sheet1, so that sheet1 ends up
# as the active/tab-selected sheet, matching the original file.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet "line_imp": insert a new row before row 4 ---
$ws2.Rows.Item(4).Insert()

# Row 3: B3/C3 turn into plain values; E3's formula is simplified to 0.05/2
$ws2.Range("B3").Value   = 0.05
$ws2.Range("C3").Value   = 0.25
$ws2.Range("E3").Formula = "=0.05/2"

# New row 4: duplicate of the (updated) row 3
$ws2.Range("A4").Value   = 23
$ws2.Range("B4").Value   = 0.05
$ws2.Range("C4").Value   = 0.25
$ws2.Range("D4").Value   = 0
$ws2.Range("E4").Formula = "=0.05/2"

# Row 5 (was row 4 before the insert): bus number corrected to 14
$ws2.Range("A5").Value = 14

# --- selections ---
$ws2.Range("D4").Select()
$ws1.Range("E9").Select()
